$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price/Volume value in this sheet is plain text (dots used as thousands
# separators, e.g. "68.662.63", or padded percentages, e.g. "  +2.32%  "). Most
# new values still read as text naturally; the handful that look like a plain
# decimal number (e.g. "0.999") need the cell pre-formatted as Text so Excel
# does not silently reinterpret the literal as a number.

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "68.708.08"
$ws.Range("D3").Value = "3.755.84"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.59"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.10"
$ws.Range("D7").Value = "3.756.32"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.31"
$ws.Range("D15").Value = "4.379.68"
$ws.Range("D16").Value = "3.755.70"
$ws.Range("D17").Value = "68.721.82"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.91"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.28"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000154"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.99"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.94"
$ws.Range("D34").Value = "3.898.58"
$ws.Range("D35").Value = "3.687.46"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.86"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "440.95"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.48"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.44"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.58"
$ws.Range("D50").Value = "2.795.85"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("E21").Value = "  +20.17%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +9.77%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("E30").Value = "  +7.86%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +2.92%  "

# --- Rows 43/44: dogwifhat moves up to rank 43, OKB drops to rank 44 ---
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.94"
$ws.Range("E43").Value = "  +6.32%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.90"
$ws.Range("E44").Value = "  +0.67%  "
